$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase (column D) values ---
# Depuration groups set first, then Uptake groups, so that sharedStrings
# ordering matches the target (Depuration=index12, Uptake=index13).
$ws.Range("D10:D27").Value = "Depuration"
$ws.Range("D33:D48").Value = "Depuration"
$ws.Range("D2:D9").Value = "Uptake"
$ws.Range("D28:D32").Value = "Uptake"

# --- Date (column E) corrections ---
$ws.Range("E36:E37").Value = 43364
$ws.Range("E41:E42").Value = 43371
$ws.Range("E44:E45").Value = 43385
$ws.Range("E46:E47").Value = 43415

# --- day-offset formulas (column F) for the lower (Figure 2B) block: anchor on
# its own first date ($E$28) instead of the upper block's $E$2. ---
$ws.Range("F3:F27").Formula = '=E3-$E$2'
$ws.Range("F28").Formula = '=E28-$E$28'
$ws.Range("F29:F48").Formula = '=E29-$E$28'

# --- toxicity (column G) value corrections, rows 2-48 ---
$ws.Range("G2").Value = 393.93939393939399
$ws.Range("G3").Value = 696.969696969697
$ws.Range("G4").Value = 939.39393939393995
$ws.Range("G5").Value = 1121.2121212121201
$ws.Range("G6").Value = 969.69696969696997
$ws.Range("G7").Value = 1424.2424242424199
$ws.Range("G8").Value = 1696.9696969696899
$ws.Range("G9").Value = 2333.3333333333298
$ws.Range("G10").Value = 2757.5757575757498
$ws.Range("G11").Value = 3878.7878787878699
$ws.Range("G12").Value = 4454.5454545454504
$ws.Range("G13").Value = 6424.2424242424204
$ws.Range("G14").Value = 8060.6060606060601
$ws.Range("G15").Value = 3636.3636363636301
$ws.Range("G16").Value = 3575.7575757575701
$ws.Range("G17").Value = 2575.7575757575701
$ws.Range("G18").Value = 2121.2121212121201
$ws.Range("G19").Value = 1787.87878787878
$ws.Range("G20").Value = 727.27272727272702
$ws.Range("G21").Value = 666.66666666666697
$ws.Range("G22").Value = 515.15151515151501
$ws.Range("G23").Value = 424.24242424242499
$ws.Range("G24").Value = 666.66666666666697
$ws.Range("G25").Value = 424.24242424242499
$ws.Range("G26").Value = 303.030303030303
$ws.Range("G27").Value = 181.81818181818099
$ws.Range("G28").Value = -0.63380281690137896
$ws.Range("G29").Value = 13.3098591549295
$ws.Range("G30").Value = 14.577464788732399
$ws.Range("G31").Value = 21.549295774647799
$ws.Range("G32").Value = 27.2535211267605
$ws.Range("G33").Value = 131.19718309859101
$ws.Range("G34").Value = 148.309859154929
$ws.Range("G35").Value = 160.98591549295699
$ws.Range("G36").Value = 156.54929577464699
$ws.Range("G37").Value = 114.718309859154
$ws.Range("G38").Value = 86.197183098591495
$ws.Range("G39").Value = 72.887323943661897
$ws.Range("G40").Value = 33.591549295774598
$ws.Range("G41").Value = 41.197183098591502
$ws.Range("G42").Value = 38.661971830985898
$ws.Range("G43").Value = 25.352112676056301
$ws.Range("G44").Value = 34.8591549295774
$ws.Range("G45").Value = 15.211267605633701
$ws.Range("G46").Value = 20.281690140845001
$ws.Range("G47").Value = 16.478873239436599
$ws.Range("G48").Value = 0.63380281690137896

# --- remove the 3rd Depuration/Figure-2B replicate (old rows 49-50): clear the
# data cells but keep E/F present (blank, original styles) ---
$ws.Range("A49:C50").ClearContents()
$ws.Range("G49:G50").ClearContents()
$ws.Range("E49:F50").ClearContents()

# --- stray styled column-I cells that appear alongside the now-blank rows 49-53 ---
$ws.Range("E49:E53").Copy()
$ws.Range("I49:I53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- selection state ---
$ws.Range("C40").Select()

Write-Output "edit complete"
